$wb = $excel.ActiveWorkbook
$s1 = $wb.Worksheets.Item("SheetGradingOrder")
$s2 = $wb.Worksheets.Item("Zipcode")
$s3 = $wb.Worksheets.Item("Zipcode_CheckOrder")

# --- SheetGradingOrder ---
$s1.Range("B1").Value = "sheet"

# --- Zipcode ---
$s2.Range("E4").Value = 43215
$s2.Range("E5").Value = 10001
$s2.Range("B6").Value = 94045
$s2.Range("E6").Value = 94043

# --- Zipcode_CheckOrder ---
$s3.Range("C1").Value = "Description"
$s3.Range("D1").Value = "special"
$s3.Range("E1").Value = "Feedback"
$s3.Range("A1").Value = "ID"

$s3.Range("D2").ClearContents()
$s3.Range("E2").ClearContents()

$s3.Range("D3").Value = "k"

$s3.Range("E4").Value = "This tab did not pass the pre-requisites. Please contact your instructor."
$s3.Range("E3").Value = "Cell B5 must be correct before this tab can be graded."
$s3.Range("F4").ClearContents()

# --- View state: selections ---
$s1.Activate() | Out-Null
$s1.Range("C5").Select() | Out-Null

$s2.Activate() | Out-Null
$s2.Range("B29").Select() | Out-Null

$s3.Activate() | Out-Null
$s3.Range("C2").Select() | Out-Null
$excel.ActiveWindow.Zoom = 94

Write-Host "done"
